$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before FC. This shifts the existing
# FC ("nom") and FD ("url_produit") columns one place to the right,
# becoming FD and FE respectively, and creates a blank FC column.
$ws.Columns("FC:FC").Insert()

# Populate the newly inserted FC column.
# Row 1 holds a new timestamp header (matching the pattern of the other
# date/time columns on that row).
$ws.Range("FC1").Value = "2026-02-04 07:37:31"

# For the data rows (2-207) the new FC column duplicates the price
# already recorded in column FB for that row (this new snapshot column
# was captured at the same price), leaving it blank where FB is blank.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $fb = $ws.Cells.Item($r, 158)   # column FB = 158
    $val = $fb.Value()
    if ($val -ne $null -and $val -ne "") {
        $ws.Cells.Item($r, 159).Value = $val   # column FC = 159
    }
}
